# Rename the two "CDC" lender options in the "Data" sheet's Prêteurs table
# (CDC froncière -> CDC pour le foncier, CDC locative -> CDC pour le logement)
# so that number-of-logements defaults to 0 instead of None.
#
# The "Data" sheet is protected, so each target cell must be temporarily
# unlocked before it can be edited, then re-locked afterwards.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

$dataSheet.Range("A5").Locked = $false
$dataSheet.Range("A5").Value2 = "CDC pour le foncier"
$dataSheet.Range("A5").Locked = $true

$dataSheet.Range("A6").Locked = $false
$dataSheet.Range("A6").Value2 = "CDC pour le logement"
$dataSheet.Range("A6").Locked = $true

# Move the selection on the Data sheet to A2 (matches the saved workbook state)
$dataSheet.Range("A2").Select() | Out-Null

# Restore "Prêts" as the active sheet/selection, as it was before our edits
$pretsSheet = $wb.Worksheets.Item("Prêts")
$pretsSheet.Activate() | Out-Null
$pretsSheet.Range("A3").Select() | Out-Null
